$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.071.92"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -2.39%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.865.91"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -2.08%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'306.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.91%  "
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'  +0.16%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.5122"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -1.79%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3749"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -0.78%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.07160"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -1.10%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.8899"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -1.19%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'20.70"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -2.80%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.07580"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.78%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.869.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -1.70%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'5.312"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -2.52%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'89.55"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -2.68%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  +0.15%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.000008456"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -2.83%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'  -2.60%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D20").Value = "'27.106.60"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -2.35%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'5.039"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -2.01%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'2.084.77"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -2.51%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D24").Value = "'6.458"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -1.88%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'1.845"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.83%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'147.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -3.50%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'17.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -1.75%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.115"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -2.46%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'112.87"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -1.36%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'4.667"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -4.04%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'4.708"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -3.04%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.09099"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +1.31%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'0.05137"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -2.87%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'3.058"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -3.70%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.157"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -6.04%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.7272"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -5.92%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.02040"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -1.99%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'2.494"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -5.40%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'3.045"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -0.90%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -1.55%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.5337"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -3.19%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'6.575"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -1.36%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'117.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.20%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'8.281"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -2.72%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.1472"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'0.4639"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -3.55%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'1.000"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.18%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'10.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -4.19%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'1.571"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -2.76%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'36.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -0.57%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'63.94"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -4.17%  "
$ws.Range("E51").ClearFormats()
